function Find-ParaIndexContaining($doc, $snippet) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($snippet)) {
            return $i
        }
    }
    return -1
}

function Add-ItalicParagraphAfter($doc, $afterIndex, $text) {
    $p = $doc.Paragraphs.Item($afterIndex)
    $p.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $newPara = $doc.Paragraphs.Item($newIndex)
    $newPara.Range.Text = $text
    $r = $doc.Range($newPara.Range.Start, $newPara.Range.End - 1)
    $r.Font.Italic = $true
    return $newIndex
}

$d = $word.ActiveDocument

# 1. Update activation date 2012 -> 2022
$d.Content.Find.Execute('Ativação: 01/01/2012', $true, $false, $false, $false, $false, $true, 1, $false, 'Ativação: 01/01/2022', 2) | Out-Null

# 2. "Objetivos" section paragraph: replace PT text, then add new italic EN paragraph after it
$idxObjetivos = Find-ParaIndexContaining $d "O TG tem por objetivo"
$d.Paragraphs.Item($idxObjetivos).Range.Text = 'O Trabalho de Graduação tem por objetivo a integração, o aprofundamento e a aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de atividades, que fazem parte do perfil de atuação profissional do engenheiro ambiental. O Trabalho de Graduação I poderá ser desenvolvido seguindo dois modelos: Modelo Artigo ou Modelo Produto.No Modelo Artigo, o aluno deverá ser capaz de apresentar a ideia principal, especificamente delimitar a questão científica, e que deve ser suportada pelo estado da arte. O projeto precisa ainda indicar o método a ser empregado na investigação com etapas previstas e os resultados esperados.No Modelo Produto, o aluno deverá ser capaz de utilizar tecnologias digitais ou outras metodologias, para obtenção de produto aplicável ou tecnologia que utilize conteúdo da Engenharia Ambiental ou áreas correlatas.'
Add-ItalicParagraphAfter $d $idxObjetivos 'The Graduate Work aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student''s ability to carry out activities, which are part of the professional performance profile of the environmental engineer. Graduation Work I can be developed following two models: Article Model or Product Model.In the Article Model, the student must be able to present the main idea, specifically to delimit the scientific question, and which must be supported by the state of the art. The project still needs to indicate the method to be used in the investigation with planned steps and the expected results.In the Product Model, the student must be able to use digital technologies or other methodologies, to obtain an applicable product or technology that uses Environmental Engineering content or related areas.' | Out-Null

# 3. "Programa resumido" paragraph: replace PT text, then add new italic EN paragraph after it
$idxResumido = Find-ParaIndexContaining $d "Elaborar a proposta de Trabalho de Gradua"
$d.Paragraphs.Item($idxResumido).Range.Text = 'O aluno deverá escolher o modelo a ser seguido para o desenvolvimento do Trabalho de Graduação I, Artigo ou Produto. Diante da escolha do modelo, o aluno deverá desenvolver conteúdo científico ou produto aplicável, ambos relacionados à Engenharia Ambiental.'
Add-ItalicParagraphAfter $d $idxResumido 'The student must choose the model to be followed for the development of the Graduation Work I, Article or Product. When choosing the model, the student must develop scientific content or applicable product, both related to Environmental Engineering.' | Out-Null

# 4. "Programa" paragraph: replace PT text, then add new italic EN paragraph after it
$idxPrograma = Find-ParaIndexContaining $d "Elaborar proposta de monografia"
$d.Paragraphs.Item($idxPrograma).Range.Text = 'Para ambos os modelos, o aluno deverá elaborar projeto cujo tema seja relacionado ao conteúdo programático do curso de Engenharia ambiental, podendo ser um tópico de interesse técnico ou científico ou ainda um estudo de caso. O programa da disciplina é constituído pelas seguintes etapas: 1. Definição e registro do tema 2. Aprovação do tema e do orientador 3. Acompanhamento mensal da evolução das atividades do aluno pelo orientador. 4. Entrega do projeto, 5. Avaliação e atribuição de notas por comissão de avaliadores.'
Add-ItalicParagraphAfter $d $idxPrograma 'For both models, the student must prepare a project whose theme is related to the programmatic content of the Environmental Engineering course, which may be a topic of technical or scientific interest or even a case study. The course program consists of the following stages: 1. Definition and registration of the theme 2. Approval of the theme and the advisor 3. Monthly monitoring of the evolution of the student''s activities by the advisor. 4. Project delivery, 5. Evaluation and grading by an evaluation committee.' | Out-Null

# 5. "Método:" run text replacement
$d.Content.Find.Execute('Em função da natureza deste curso, a avaliação será feita pela elaboração e apresentação de um plano de trabalho.', $true, $false, $false, $false, $false, $true, 1, $false, 'Modelo Artigo: O aluno deverá apresentar um projeto de pesquisa seguindo modelo fornecido pela responsável da disciplina. A avaliação do projeto será realizada utilizando o sistema online da Escola de Engenharia de Lorena, por dois avaliadores doutores ou especialistas na área do projeto com elaboração de parecer. Os avaliadores deverão analisar o referencial teórico, exequibilidade do projeto e a aderência ao curso de Engenharia Ambiental em comparação com a ideia principal e a questão científica indicados pelo aluno.Modelo Produto: O aluno deverá apresentar um projeto de desenvolvimento de produto: aplicativo Móvel; Livro Digital; Animação; Jogos Didáticos, ou qualquer outro produto relacionado com a Engenharia Ambiental ou área correlata. Para apresentação do projeto deverá ser utilizado modelo fornecido pela responsável da disciplina. A avaliação do projeto será realizada utilizando o sistema online da Escola de Engenharia de Lorena, por dois avaliadores doutores ou especialistas na área do projeto com elaboração de parecer. Os avaliadores deverão analisar o referencial teórico, exequibilidade do projeto, uso de tecnologia no desenvolvimento e a aderência ao curso de Engenharia Ambiental.Em ambos modelos, a responsável pela disciplina irá emitir uma terceira nota, em que o aluno será avaliado pelo seu desempenho na disciplina durante o semestre.', 2) | Out-Null

# 6. "Critério:" run text replacement
$d.Content.Find.Execute('Avaliação e atribuição de nota do Trabalho de Graduação por uma comissão de professores.', $true, $false, $false, $false, $false, $true, 1, $false, 'Avaliação, elaboração de parecer e atribuição de nota por dois avaliadores doutores ou especialistas na área do projeto, bem como pela responsável da disciplina. A nota final será a média das três notas de igual pesoFica sob responsabilidade do orientador a verificação de ocorrência de plágio utilizando software apropriado e avaliação em Comitê de Ética, quando exigido, via cadastro na Plataforma Brasil.', 2) | Out-Null
